$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells I1 and J1 - copy H1's format (bold/border/center-top) then set text
$ws.Range("H1").Copy($ws.Range("I1"))
$ws.Range("H1").Copy($ws.Range("J1"))
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for I (I0) and J (IF) columns, rows 2-27, in row order
$rows = 2..27
$iVals = @(5, 4, 9, 3, 1, 2, 9, 6, 6, 6, 7, 3, 5, 7, 7, 6, 6, 6, 9, 6, 8, 8, 7, 4, 5, 3)
$jVals = @(6, 5, 9, 3, 1, 2, 9, 6, 6, 6, 7, 4, 5, 7, 7, 6, 6, 6, 9, 6, 8, 8, 8, 4, 5, 3)

for ($idx = 0; $idx -lt $rows.Count; $idx++) {
    $row = $rows[$idx]
    $ws.Cells.Item($row, 9).Value = $iVals[$idx]
    $ws.Cells.Item($row, 10).Value = $jVals[$idx]
}
